# Actualización automática 2025-11-18 16:30:08
# Inserts a new advisor "LINDAO RODRIGUEZ ANTONIO COLON" (all-zero sales) in both
# "VENTAS POR GRUPO" and "VENTA MENSUAL" sheets right before "LOAIZA TINOCO JUAN
# PABLO" (row 33), shifting the remaining advisors + totals row down by one, and
# refreshes a handful of figures that changed for this reporting run.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "VENTAS POR GRUPO"
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")

# Figures updated for ECUAFERRI S.A. (row 21)
$ws1.Range("H21").Value = 423.9
$ws1.Range("I21").Value = 503.33

# Insert the new advisor row before the current row 33 (LOAIZA TINOCO JUAN PABLO).
# Excel shifts LOAIZA..ZAVALA down one row (33->34 ... 54->55) and the trailing
# "X de 53" summary row moves from 55 to 56.
$ws1.Rows.Item(33).Insert()

$ws1.Range("A33").Value = "OFICINA-CATAECSA"
$ws1.Range("B33").Value = "LINDAO RODRIGUEZ ANTONIO COLON"
$ws1.Range("C33").Value = 0
$ws1.Range("D33").Value = 0
$ws1.Range("E33").Value = 0
$ws1.Range("F33").Value = 0
$ws1.Range("G33").Value = 0
$ws1.Range("H33").Value = 0
$ws1.Range("I33").Value = 0
$ws1.Range("J33").Value = 0
$ws1.Range("K33").Value = 0
$ws1.Range("L33").Value = 0
$ws1.Range("M33").Value = 0
$ws1.Range("N33").Value = 0
$ws1.Range("O33").Value = 0
$ws1.Range("P33").Value = 0
$ws1.Range("Q33").Value = 0
$ws1.Range("R33").Value = 0

# Refresh the "X de 53" -> "X de 54" counters on the (now) last row, 56.
$ws1.Range("C56").Value = "0 de 54"
$ws1.Range("D56").Value = "0 de 54"
$ws1.Range("E56").Value = "1 de 54"
$ws1.Range("F56").Value = "0 de 54"
$ws1.Range("G56").Value = "0 de 54"
$ws1.Range("H56").Value = "1 de 54"
$ws1.Range("I56").Value = "1 de 54"
$ws1.Range("J56").Value = "0 de 54"
$ws1.Range("K56").Value = "0 de 54"
$ws1.Range("L56").Value = "1 de 54"
$ws1.Range("M56").Value = "1 de 54"
$ws1.Range("N56").Value = "0 de 54"
$ws1.Range("O56").Value = "0 de 54"
$ws1.Range("P56").Value = "0 de 54"
$ws1.Range("Q56").Value = "0 de 54"
$ws1.Range("R56").Value = "0 de 54"

# ---------------------------------------------------------------------------
# Sheet "VENTA MENSUAL"
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")

# Figure updated for ECUAFERRI S.A. (row 21)
$ws2.Range("F21").Value = 927.23

# Same insertion as above.
$ws2.Rows.Item(33).Insert()

$ws2.Range("A33").Value = "OFICINA-CATAECSA"
$ws2.Range("B33").Value = "LINDAO RODRIGUEZ ANTONIO COLON"
$ws2.Range("C33").Value = 0
$ws2.Range("D33").Value = 0
$ws2.Range("E33").Value = 0
$ws2.Range("F33").Value = 0
$ws2.Range("G33").Value = 0

# Totals row moved from 55 to 56; only the "noviembre" total changed.
$ws2.Range("F56").Value = 2177.85

# ---------------------------------------------------------------------------
# Sheet "CUMPLIMIENTO MENSUAL"
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

$ws3.Range("D3").Value = 4518.7
$ws3.Range("E3").Value = -4518.7
$ws3.Range("D5").Value = 13296.88
$ws3.Range("E5").Value = 12789.53
$ws3.Range("F5").Value = 0.509724412059766
